# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" (so it becomes the
#    second tab, pushing 2022-Q2 / 2022-Q1 / ... down by one position).
# 2) Populate "总计" with a new summary row for 2022-Q4 at the top of the
#    data (row 2), shifting the existing quarters down by one row.
# 3) Fill the new "2022-Q4" sheet with its fund-holding detail rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new sheet right after "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet - insert a row for 2022-Q4
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# restore formatting for the new row from the row below (old row2, now row3)
$total.Range("A3:D3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# The row-index column (A) is a simple 0-based running counter; re-number it
# for every data row (the shift-down from Insert() only moved old values
# down with the rows, it did not renumber them).
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(8,1).Value = 6

$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.04

# ---------------------------------------------------------------------
# Step 3: populate the new "2022-Q4" sheet
# ---------------------------------------------------------------------
$donor = $wb.Worksheets.Item("2022-Q2")
$q4 = $wb.Worksheets.Item("2022-Q4")

# Copy header/row formatting (styles + borders) from the donor sheet, which
# has the identical table layout. (Column A on the header row, and columns
# B..H on the index rows, carry no special formatting in the source sheets,
# so copy the two formatted blocks separately to avoid manufacturing empty
# placeholder cells that don't exist in the donor.)
$donor.Range("B1:H3").Copy()
$q4.Range("B1").PasteSpecial(-4122)
$donor.Range("A2:A3").Copy()
$q4.Range("A2").PasteSpecial(-4122)

$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Fund code + numeric-looking columns must stay text (leading zeros / fixed
# decimal formatting), so force text format before writing them.
$q4.Range("B2:B3").NumberFormat = "@"
$q4.Range("D2:G3").NumberFormat = "@"

$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "006336"
$q4.Cells.Item(2,3).Value = "泓德量化精选混合"
$q4.Cells.Item(2,4).Value = "2.55"
$q4.Cells.Item(2,5).Value = "93.68"
$q4.Cells.Item(2,6).Value = "1.56"
$q4.Cells.Item(2,7).Value = "0.0398"
$q4.Cells.Item(2,8).Value = 8

$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "080015"
$q4.Cells.Item(3,3).Value = "长盛中小盘精选混合"
$q4.Cells.Item(3,4).Value = "0.13"
$q4.Cells.Item(3,5).Value = "84.17"
$q4.Cells.Item(3,6).Value = "2.35"
$q4.Cells.Item(3,7).Value = "0.0031"
$q4.Cells.Item(3,8).Value = 8
